# Update Financials - Yearly RELV: refresh figures pulled from source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - Cost of Revenue
$ws.Range("D9").Value = 24100
$ws.Range("E9").Value = 26100
$ws.Range("F9").Value = 29500
$ws.Range("G9").Value = 32200
$ws.Range("H9").Value = 38900
$ws.Range("I9").Value = 39500
$ws.Range("J9").Value = 42700

# Row 10 - Gross Profit
$ws.Range("D10").Value = 17700
$ws.Range("E10").Value = 19400
$ws.Range("F10").Value = 22300
$ws.Range("G10").Value = 25100
$ws.Range("H10").Value = 29300
$ws.Range("I10").Value = 29200
$ws.Range("J10").Value = 31100

# Row 83 - Depreciation
$ws.Range("D83").Value = 1000
$ws.Range("E83").Value = 1000
$ws.Range("F83").Value = 1000
$ws.Range("G83").Value = 1000
$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 1100
$ws.Range("J83").Value = "NA"

# Row 89 - Total Cash Flow From Operating Activities
$ws.Range("D89").Value = 1500
$ws.Range("E89").Value = -800
$ws.Range("F89").Value = -400
$ws.Range("G89").Value = 2500
$ws.Range("H89").Value = 2500
$ws.Range("I89").Value = 2800
$ws.Range("J89").Value = 2200

# Row 91 - Capital Expenditures
$ws.Range("D91").Value = -200
$ws.Range("E91").Value = -300
$ws.Range("F91").Value = -900
$ws.Range("G91").Value = -400
$ws.Range("H91").Value = -500
$ws.Range("I91").Value = -400
$ws.Range("J91").Value = -600

# Row 94 - Total Cash Flows From Investing Activities
$ws.Range("D94").Value = -100
$ws.Range("E94").Value = -100
$ws.Range("F94").Value = -1100
$ws.Range("G94").Value = -1700
$ws.Range("H94").Value = -2800
$ws.Range("I94").Value = -700
$ws.Range("J94").Value = "NA"

# Row 96 - Dividends Paid
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = -400
$ws.Range("H96").Value = -400
$ws.Range("I96").Value = -500
$ws.Range("J96").Value = -500

# Row 100 - Total Cash Flows From Financing Activities
$ws.Range("D100").Value = -1000
$ws.Range("E100").Value = -700
$ws.Range("F100").Value = -100
$ws.Range("G100").Value = 100
$ws.Range("H100").Value = -1200
$ws.Range("I100").Value = -1200
$ws.Range("J100").Value = "NA"

# Row 101 - Effect Of Exchange Rate Changes
$ws.Range("D101").Value = -100
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = -100
$ws.Range("G101").Value = -100
$ws.Range("H101").Value = 100
$ws.Range("I101").Value = -100
$ws.Range("J101").Value = "NA"

# Row 102 - Change In Cash and Cash Equivalents
$ws.Range("D102").Value = 300
$ws.Range("E102").Value = -1700
$ws.Range("F102").Value = -1700
$ws.Range("G102").Value = 900
$ws.Range("H102").Value = -1400
$ws.Range("I102").Value = 800
$ws.Range("J102").Value = 600
